$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "52"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "56"
$find.Execute(
    "52",    # FindText
    $true,   # MatchCase
    $true,   # MatchWholeWord
    $false,  # MatchWildcards
    $false,  # MatchSoundsLike
    $false,  # MatchAllWordForms
    $true,   # Forward
    1,       # Wrap
    $false,  # Format
    "56",    # ReplaceWith
    2        # Replace (wdReplaceAll)
)
